# Insert two new rows at row 123 (this shifts the existing rows 123-151 down to 125-153,
# preserving all their original values and formatting).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A123:A124").EntireRow.Insert()

# Fill in the new row 123 with its data (Primera quality, new price observation).
$ws.Range("A123").Value = 2
$ws.Range("B123").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C123").Value = "Coquimbo"
$ws.Range("D123").Value = 44637
$ws.Range("E123").Value = 4
$ws.Range("F123").Value = 100112043
$ws.Range("G123").Value = "Pepino ensalada"
$ws.Range("H123").Value = "Sin especificar"
$ws.Range("I123").Value = "Primera"
$ws.Range("J123").Value = 240
$ws.Range("K123").Value = 16000
$ws.Range("L123").Value = 17000
$ws.Range("M123").Value = 16500
$ws.Range("N123").Value = "`$/caja 70 unidades"
$ws.Range("O123").Value = "Provincia de Limarí"
$ws.Range("P123").Value = 236
$ws.Range("Q123").Value = 70
$ws.Range("R123").Value = "Hortaliza"

# Fill in the new row 124 with its data (Segunda quality, new price observation).
$ws.Range("A124").Value = 2
$ws.Range("B124").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C124").Value = "Coquimbo"
$ws.Range("D124").Value = 44637
$ws.Range("E124").Value = 4
$ws.Range("F124").Value = 100112043
$ws.Range("G124").Value = "Pepino ensalada"
$ws.Range("H124").Value = "Sin especificar"
$ws.Range("I124").Value = "Segunda"
$ws.Range("J124").Value = 200
$ws.Range("K124").Value = 14000
$ws.Range("L124").Value = 15000
$ws.Range("M124").Value = 14500
$ws.Range("N124").Value = "`$/caja 100 unidades"
$ws.Range("O124").Value = "Provincia de Limarí"
$ws.Range("P124").Value = 145
$ws.Range("Q124").Value = 100
$ws.Range("R124").Value = "Hortaliza"
